$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row from diff hunk @@ -3757,25 +3757,25 @@
$ws.Range("H62").Value = 3292
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 3365
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 3365
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4613
# row from diff hunk @@ -3907,25 +3907,25 @@
$ws.Range("H65").Value = 3292
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 3365
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 16825
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -23065
# row from diff hunk @@ -5590,25 +5590,25 @@
$ws.Range("H98").Value = 3744.8667
$ws.Range("I98").Value = 2782.5386
$ws.Range("J98").Value = 10000
$ws.Range("K98").Value = 2782.5386
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = -1284.5386
$ws.Range("N98").Value = -12996
# row from diff hunk @@ -6747,25 +6747,25 @@
$ws.Range("H121").Value = 1250
$ws.Range("J121").Value = 2100
$ws.Range("L121").Value = 6300
$ws.Range("N121").Value = -9794
# row from diff hunk @@ -6799,25 +6799,25 @@
$ws.Range("H122").Value = 3744.8667
$ws.Range("I122").Value = 2782.5386
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 8347.6158
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -5897.6158
$ws.Range("N122").Value = -34900
# row from diff hunk @@ -7295,25 +7295,25 @@
$ws.Range("H132").Value = 10007690
$ws.Range("I132").Value = 12508275
$ws.Range("J132").Value = 5350
$ws.Range("K132").Value = 37524825
$ws.Range("L132").Value = 16050
$ws.Range("M132").Value = -37522295
$ws.Range("N132").Value = -21110
# row from diff hunk @@ -7396,22 +7396,22 @@
$ws.Range("H134").Value = 30000
$ws.Range("J134").Value = 30000
$ws.Range("L134").Value = 30000
$ws.Range("N134").Value = -40140

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row from diff hunk @@ -9930,22 +9930,19 @@
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = ""
# row from diff hunk @@ -13807,25 +13804,25 @@
$ws.Range("H122").Value = 2675.8333
$ws.Range("I122").Value = 1932.3684
$ws.Range("J122").Value = 3960
$ws.Range("K122").Value = 5797.1052
$ws.Range("L122").Value = 11880
$ws.Range("M122").Value = -3347.1052
$ws.Range("N122").Value = -16780
# row from diff hunk @@ -14294,25 +14291,25 @@
$ws.Range("H132").Value = 40003612
$ws.Range("I132").Value = 71431470
$ws.Range("J132").Value = 4518
$ws.Range("K132").Value = 214294410
$ws.Range("L132").Value = 13554
$ws.Range("M132").Value = -214291880
$ws.Range("N132").Value = -18614

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row from diff hunk @@ -23411,22 +23408,22 @@
$ws.Range("H33").Value = 42920
$ws.Range("I33").Value = 2247.5
$ws.Range("K33").Value = 2247.5
$ws.Range("M33").Value = -1868.5
# row from diff hunk @@ -28387,25 +28384,25 @@
$ws.Range("H134").Value = 1338.8975
$ws.Range("I134").Value = 838.2
$ws.Range("J134").Value = 2233
$ws.Range("K134").Value = 2514.6
$ws.Range("L134").Value = 6699
$ws.Range("M134").Value = 20.39999999999964
$ws.Range("N134").Value = -11769

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row from diff hunk @@ -29136,22 +29133,25 @@
$ws.Range("H7").Value = 80
$ws.Range("I7").Value = 76.666664
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 229.999992
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -117.999992
$ws.Range("N7").Value = -524
# row from diff hunk @@ -29644,25 +29644,22 @@
$ws.Range("H17").Value = 525
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 525
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1575
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = -1913
# row from diff hunk @@ -30057,25 +30054,25 @@
$ws.Range("H25").Value = 1266.3334
$ws.Range("I25").Value = 199
$ws.Range("J25").Value = 1800
$ws.Range("K25").Value = 597
$ws.Range("L25").Value = 5400
$ws.Range("M25").Value = -428
$ws.Range("N25").Value = -5738
# row from diff hunk @@ -30314,25 +30311,25 @@
$ws.Range("H30").Value = 1266.3334
$ws.Range("I30").Value = 199
$ws.Range("J30").Value = 1800
$ws.Range("K30").Value = 597
$ws.Range("L30").Value = 5400
$ws.Range("M30").Value = -495
$ws.Range("N30").Value = -5604
# row from diff hunk @@ -30522,25 +30519,25 @@
$ws.Range("H34").Value = 11037.9
$ws.Range("J34").Value = 12208.777
$ws.Range("L34").Value = 36626.331
$ws.Range("N34").Value = -36794.331
# row from diff hunk @@ -30776,25 +30773,22 @@
$ws.Range("H39").Value = 2380.8
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2380.8
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 7142.400000000001
$ws.Range("M39").Value = ""
$ws.Range("N39").Value = -7730.400000000001
# row from diff hunk @@ -32652,25 +32646,25 @@
$ws.Range("H76").Value = 2902.1667
$ws.Range("I76").Value = 1837.6666
$ws.Range("J76").Value = 3966.6667
$ws.Range("K76").Value = 5512.9998
$ws.Range("L76").Value = 11900.0001
$ws.Range("M76").Value = -5129.9998
$ws.Range("N76").Value = -12666.0001
# row from diff hunk @@ -32805,25 +32799,25 @@
$ws.Range("H79").Value = 2902.1667
$ws.Range("I79").Value = 1837.6666
$ws.Range("J79").Value = 3966.6667
$ws.Range("K79").Value = 5512.9998
$ws.Range("L79").Value = 11900.0001
$ws.Range("M79").Value = -4186.9998
$ws.Range("N79").Value = -14552.0001
# row from diff hunk @@ -32857,22 +32851,25 @@
$ws.Range("H80").Value = 6901
$ws.Range("I80").Value = 1800
$ws.Range("J80").Value = 9451.5
$ws.Range("K80").Value = 5400
$ws.Range("L80").Value = 28354.5
$ws.Range("M80").Value = -4464
$ws.Range("N80").Value = -30226.5
# row from diff hunk @@ -33007,22 +33004,25 @@
$ws.Range("H83").Value = 6901
$ws.Range("I83").Value = 1800
$ws.Range("J83").Value = 9451.5
$ws.Range("K83").Value = 16200
$ws.Range("L83").Value = 85063.5
$ws.Range("M83").Value = -11520
$ws.Range("N83").Value = -94423.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row from diff hunk @@ -37859,23 +37859,20 @@
$ws.Range("H36").Value = 4500
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = ""
# row from diff hunk @@ -38208,25 +38205,25 @@
$ws.Range("H43").Value = 2683.4
$ws.Range("I43").Value = 1205.6666
$ws.Range("J43").Value = 4900
$ws.Range("K43").Value = 1205.6666
$ws.Range("L43").Value = 4900
$ws.Range("M43").Value = -1054.6666
$ws.Range("N43").Value = -5202
# row from diff hunk @@ -38361,22 +38358,22 @@
$ws.Range("H46").Value = 9948.666999999999
$ws.Range("J46").Value = 9948.666999999999
$ws.Range("L46").Value = 9948.666999999999
$ws.Range("N46").Value = -10260.667
# row from diff hunk @@ -42245,22 +42242,22 @@
$ws.Range("H126").Value = 3319.8096
$ws.Range("I126").Value = 2127.625
$ws.Range("K126").Value = 6382.875
$ws.Range("M126").Value = -3912.875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row from diff hunk @@ -46038,25 +46035,25 @@
$ws.Range("H61").Value = 31251870
$ws.Range("I61").Value = 38463148
$ws.Range("J61").Value = 3000.8333
$ws.Range("K61").Value = 38463148
$ws.Range("L61").Value = 3000.8333
$ws.Range("M61").Value = -38462946
$ws.Range("N61").Value = -3404.8333
# row from diff hunk @@ -48577,25 +48574,25 @@
$ws.Range("H113").Value = 31251870
$ws.Range("I113").Value = 38463148
$ws.Range("J113").Value = 3000.8333
$ws.Range("K113").Value = 38463148
$ws.Range("L113").Value = 3000.8333
$ws.Range("M113").Value = -38460978
$ws.Range("N113").Value = -7340.8333
# row from diff hunk @@ -49009,25 +49006,25 @@
$ws.Range("H122").Value = 3103.7144
$ws.Range("I122").Value = 2550.2222
$ws.Range("J122").Value = 4100
$ws.Range("K122").Value = 7650.6666
$ws.Range("L122").Value = 12300
$ws.Range("M122").Value = -5200.6666
$ws.Range("N122").Value = -17200
# row from diff hunk @@ -49505,22 +49502,22 @@
$ws.Range("H132").Value = 3510.8462
$ws.Range("I132").Value = 2323.6667
$ws.Range("K132").Value = 6971.000100000001
$ws.Range("M132").Value = -4441.000100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row from diff hunk @@ -54490,22 +54487,19 @@
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""
# row from diff hunk @@ -56015,22 +56009,22 @@
$ws.Range("H123").Value = 21783.691
$ws.Range("J123").Value = 21783.691
$ws.Range("L123").Value = 21783.691
$ws.Range("N123").Value = -31583.691
